# Weekly-style refresh of the cryptos price/volume table, matching the
# GitHub Actions commit "Updated cryptos list ... with GitHub Actions".
#
# Column D ("Price") and column E ("Volume(1h)") are plain-text cells in
# the source sheet (prices use '.' as both a thousands AND decimal
# separator, e.g. "66.848.47", so they can never be real numbers; the
# percentages keep two leading/trailing spaces of padding). When COM
# assigns a Value that *looks* numeric (e.g. "592.19"), Excel will
# normally coerce it to a real number - so for those cells we force the
# format to Text first and restore the "Normal" style afterwards so no
# stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value
    )
    $range = $ws.Range($Cell)
    $range.NumberFormat = "@"
    $range.Value = $Value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.997.91"
$ws.Range("E2").Value = "  -0.53%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.454.43"
$ws.Range("E3").Value = "  -1.69%  "

# Row 4 - TetherUSD (price unchanged)
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "592.19"
$ws.Range("E5").Value = "  -1.17%  "

# Row 6 - Solana
Set-TextValue "D6" "178.67"
$ws.Range("E6").Value = "  +2.67%  "

# Row 7 - was USDC, now XRP
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D7" "0.611"
$ws.Range("E7").Value = "  +4.55%  "

# Row 8 - was XRP, now USDC
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - LidoStakedEther (volume unchanged)
Set-TextValue "D9" "3.455.22"

# Row 10 - Dogecoin
Set-TextValue "D10" "0.137"
$ws.Range("E10").Value = "  +3.03%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.94"
$ws.Range("E11").Value = "  -3.37%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.431"
$ws.Range("E12").Value = "  -0.24%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.052.53"
$ws.Range("E13").Value = "  -1.44%  "

# Row 14 - Avalanche
Set-TextValue "D14" "31.50"
$ws.Range("E14").Value = "  +4.99%  "

# Row 15 - TRON (price unchanged)
$ws.Range("E15").Value = "  -0.41%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "67.034.76"
$ws.Range("E16").Value = "  -0.39%  "

# Row 17 - ShibaInu (price unchanged)
$ws.Range("E17").Value = "  -2.56%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.454.86"
$ws.Range("E18").Value = "  -1.49%  "

# Row 19 - Polkadot (volume unchanged)
Set-TextValue "D19" "6.24"

# Row 20 - Chainlink
Set-TextValue "D20" "14.10"
$ws.Range("E20").Value = "  -2.69%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "387.19"
$ws.Range("E21").Value = "  -1.66%  "

# Row 22 - Uniswap
Set-TextValue "D22" "7.90"
$ws.Range("E22").Value = "  -1.03%  "

# Row 23 - Dai
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.13%  "

# Row 24 - LEO
Set-TextValue "D24" "5.76"
$ws.Range("E24").Value = "  +1.24%  "

# Row 25 - Litecoin
Set-TextValue "D25" "72.23"
$ws.Range("E25").Value = "  -1.75%  "

# Row 26 - Polygon
Set-TextValue "D26" "0.534"
$ws.Range("E26").Value = "  -0.56%  "

# Row 27 - PEPE
Set-TextValue "D27" "0.0000121"
$ws.Range("E27").Value = "  -1.13%  "

# Row 28 - InternetComputer(DFINITY) (price unchanged)
$ws.Range("E28").Value = "  +0.37%  "

# Row 29 - Kaspa (price unchanged)
$ws.Range("E29").Value = "  -3.98%  "

# Row 30 - Binance-PegBSC-USD (price unchanged)
$ws.Range("E30").Value = "  +0.40%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "6.14"
$ws.Range("E31").Value = "  -0.31%  "

# Row 32 - PancakeSwap (price unchanged)
$ws.Range("E32").Value = "  -0.73%  "

# Row 33 - Fetch.AI (price unchanged)
$ws.Range("E33").Value = "  -2.80%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "23.37"
$ws.Range("E34").Value = "  -1.35%  "

# Row 35 - Aptos
Set-TextValue "D35" "7.32"
$ws.Range("E35").Value = "  -1.14%  "

# Row 36 - USDe: unchanged, nothing to do.

# Row 37 - ImmutableX (price unchanged)
$ws.Range("E37").Value = "  -2.46%  "

# Row 38 - Monero
Set-TextValue "D38" "161.86"
$ws.Range("E38").Value = "  -0.99%  "

# Row 39 - Mantle
Set-TextValue "D39" "0.875"
$ws.Range("E39").Value = "  -0.37%  "

# Row 40 - dogwifhat (price unchanged)
$ws.Range("E40").Value = "  +7.66%  "

# Row 41 - Stacks (price unchanged)
$ws.Range("E41").Value = "  -3.98%  "

# Row 42 - RenderToken
Set-TextValue "D42" "6.80"
$ws.Range("E42").Value = "  -1.72%  "

# Row 43 - Filecoin
Set-TextValue "D43" "4.64"
$ws.Range("E43").Value = "  -0.63%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "26.06"
$ws.Range("E44").Value = "  -0.28%  "

# Row 45 - was Maker, now Hedera
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D45" "0.0718"
$ws.Range("E45").Value = "  -1.86%  "

# Row 46 - was Hedera, now Maker
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D46" "2.756.74"
$ws.Range("E46").Value = "  -2.20%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "25.90"
$ws.Range("E47").Value = "  -4.90%  "

# Row 48 - OKB
Set-TextValue "D48" "41.11"
$ws.Range("E48").Value = "  -3.38%  "

# Row 49 - VeChain
Set-TextValue "D49" "0.0297"
$ws.Range("E49").Value = "  -2.62%  "

# Row 50 - Bittensor
Set-TextValue "D50" "325.85"
$ws.Range("E50").Value = "  -5.18%  "

# Row 51 - ONDO (price unchanged)
$ws.Range("E51").Value = "  -4.00%  "
